# New crime data collected — weekly CompStat (Bronx) update.
# Updates: report header (issue number + week-covering dates), and the
# Week-to-Date / 28-Day / Year-to-Date / 2-Year crime-complaint figures
# in the main table (rows 14-33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (rich-text cells): only the embedded numbers/dates change,
# the surrounding label text ("Volume 32   Number  ", "Report Covering
# the Week  " / "  Through  ") stays the same. Use Characters() to swap
# just the affected substring in place rather than replacing the whole
# cell value.
# ---------------------------------------------------------------------

# A8: "Volume 32   Number  20" -> "Volume 32   Number  21"
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "21"

# C9: "Report Covering the Week  5/12/2025  Through  5/18/2025"
#  -> "Report Covering the Week  5/19/2025  Through  5/25/2025"
$c9 = $ws.Range("C9")
$c9.Characters(27, 9).Text = "5/19/2025"
$c9.Characters(47, 9).Text = "5/25/2025"

# ---------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = -80
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = 13
$ws.Range("H14").Value = -38.461538461538
$ws.Range("I14").Value = 37
$ws.Range("J14").Value = 45
$ws.Range("K14").Value = -17.777777777777
$ws.Range("L14").Value = -26
$ws.Range("M14").Value = -13.953488372093
$ws.Range("N14").Value = -79.781420765027

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
$ws.Range("C15").Value = 9
$ws.Range("D15").Value = 10
$ws.Range("E15").Value = -10
$ws.Range("F15").Value = 44
$ws.Range("G15").Value = 36
$ws.Range("H15").Value = 22.222222222222
$ws.Range("I15").Value = 206
$ws.Range("J15").Value = 167
$ws.Range("K15").Value = 23.353293413173
$ws.Range("L15").Value = 26.380368098159
$ws.Range("M15").Value = 85.585585585585
$ws.Range("N15").Value = -23.420074349442

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 97
$ws.Range("D16").Value = 93
$ws.Range("E16").Value = 4.301075268817
$ws.Range("F16").Value = 420
$ws.Range("G16").Value = 367
$ws.Range("H16").Value = 14.441416893733
$ws.Range("I16").Value = 1759
$ws.Range("J16").Value = 1893
$ws.Range("K16").Value = -7.078711040676
$ws.Range("L16").Value = 1.911935110081
$ws.Range("M16").Value = 7.321537522879
$ws.Range("N16").Value = -72.189723320158

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 155
$ws.Range("D17").Value = 211
$ws.Range("E17").Value = -26.540284360189
$ws.Range("F17").Value = 677
$ws.Range("G17").Value = 688
$ws.Range("H17").Value = -1.598837209302
$ws.Range("I17").Value = 3217
$ws.Range("J17").Value = 3069
$ws.Range("K17").Value = 4.822417725643
$ws.Range("L17").Value = 8.280040390440
$ws.Range("M17").Value = 90.580568720379
$ws.Range("N17").Value = -2.010356381358

# ---------------------------------------------------------------------
# Row 18 - Burglary  (D18 unchanged = 64)
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 49
$ws.Range("E18").Value = -23.4375
$ws.Range("F18").Value = 182
$ws.Range("G18").Value = 226
$ws.Range("H18").Value = -19.469026548672
$ws.Range("I18").Value = 1119
$ws.Range("J18").Value = 1143
$ws.Range("K18").Value = -2.099737532808
$ws.Range("L18").Value = -8.727569331158
$ws.Range("M18").Value = -6.594323873121
$ws.Range("N18").Value = -84.845612134344

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 177
$ws.Range("D19").Value = 196
$ws.Range("E19").Value = -9.693877551020
$ws.Range("F19").Value = 735
$ws.Range("G19").Value = 719
$ws.Range("H19").Value = 2.225312934631
$ws.Range("I19").Value = 3489
$ws.Range("J19").Value = 3562
$ws.Range("K19").Value = -2.049410443571
$ws.Range("L19").Value = 14.543663821405
$ws.Range("M19").Value = 108.672248803828
$ws.Range("N19").Value = 26.230101302460

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 118
$ws.Range("D20").Value = 81
$ws.Range("E20").Value = 45.679012345679
$ws.Range("F20").Value = 420
$ws.Range("G20").Value = 311
$ws.Range("H20").Value = 35.048231511254
$ws.Range("I20").Value = 1754
$ws.Range("J20").Value = 1575
$ws.Range("K20").Value = 11.365079365079
$ws.Range("L20").Value = -17.652582159624
$ws.Range("M20").Value = 123.724489795918
$ws.Range("N20").Value = -71.278860324218

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 606
$ws.Range("D21").Value = 660
$ws.Range("E21").Value = -8.181818181818
$ws.Range("F21").Value = 2486
$ws.Range("G21").Value = 2360
$ws.Range("H21").Value = 5.338983050847
$ws.Range("I21").Value = 11581
$ws.Range("J21").Value = 11454
$ws.Range("K21").Value = 1.108782957918
$ws.Range("L21").Value = 2.378005657708
$ws.Range("M21").Value = 62.312543798178
$ws.Range("N21").Value = -55.990879726391

# ---------------------------------------------------------------------
# Row 22 - Transit  (F22 unchanged = 21, N22 unchanged "0")
# ---------------------------------------------------------------------
$ws.Range("C22").Value = 6
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = -14.285714285714
$ws.Range("G22").Value = 22
$ws.Range("H22").Value = -4.545454545454
$ws.Range("I22").Value = 117
$ws.Range("J22").Value = 142
$ws.Range("K22").Value = -17.605633802816
$ws.Range("L22").Value = -0.847457627118
$ws.Range("M22").Value = -10.687022900763

# ---------------------------------------------------------------------
# Row 23 - Housing  (N23 unchanged "0")
# ---------------------------------------------------------------------
$ws.Range("C23").Value = 24
$ws.Range("D23").Value = 33
$ws.Range("E23").Value = -27.272727272727
$ws.Range("F23").Value = 115
$ws.Range("G23").Value = 115
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 597
$ws.Range("J23").Value = 668
$ws.Range("K23").Value = -10.628742514970
$ws.Range("L23").Value = -14.592274678111
$ws.Range("M23").Value = 56.282722513089

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny  (N24 unchanged "0")
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 398
$ws.Range("D24").Value = 287
$ws.Range("E24").Value = 38.675958188153
$ws.Range("F24").Value = 1397
$ws.Range("G24").Value = 1206
$ws.Range("H24").Value = 15.837479270315
$ws.Range("I24").Value = 6916
$ws.Range("J24").Value = 6534
$ws.Range("K24").Value = 5.846342209978
$ws.Range("L24").Value = -0.417566594672
$ws.Range("M24").Value = 46.184738955823

# ---------------------------------------------------------------------
# Row 25 - Retail Theft  (M25, N25 unchanged "0")
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 116
$ws.Range("D25").Value = 91
$ws.Range("E25").Value = 27.472527472527
$ws.Range("F25").Value = 416
$ws.Range("G25").Value = 457
$ws.Range("H25").Value = -8.971553610503
$ws.Range("I25").Value = 2259
$ws.Range("J25").Value = 2626
$ws.Range("K25").Value = -13.975628332064
$ws.Range("L25").Value = -24.900265957446

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault  (N26 unchanged "0")
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 249
$ws.Range("D26").Value = 286
$ws.Range("E26").Value = -12.937062937062
$ws.Range("F26").Value = 979
$ws.Range("G26").Value = 926
$ws.Range("H26").Value = 5.723542116630
$ws.Range("I26").Value = 4229
$ws.Range("J26").Value = 4204
$ws.Range("K26").Value = 0.594671741198
$ws.Range("L26").Value = 4.471343873517
$ws.Range("M26").Value = 1.172248803827

# ---------------------------------------------------------------------
# Row 27 - UCR Rape*  (M27, N27 unchanged "0")
# ---------------------------------------------------------------------
$ws.Range("C27").Value = 11
$ws.Range("D27").Value = 14
$ws.Range("E27").Value = -21.428571428571
$ws.Range("F27").Value = 52
$ws.Range("G27").Value = 54
$ws.Range("H27").Value = -3.703703703703
$ws.Range("I27").Value = 260
$ws.Range("J27").Value = 265
$ws.Range("K27").Value = -1.886792452830
$ws.Range("L27").Value = -2.255639097744

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes  (M28, N28 unchanged "0")
# ---------------------------------------------------------------------
$ws.Range("C28").Value = 24
$ws.Range("D28").Value = 20
$ws.Range("E28").Value = 20
$ws.Range("F28").Value = 99
$ws.Range("G28").Value = 107
$ws.Range("H28").Value = -7.476635514018
$ws.Range("I28").Value = 423
$ws.Range("J28").Value = 474
$ws.Range("K28").Value = -10.759493670886
$ws.Range("L28").Value = 0.954653937947

# ---------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 10
$ws.Range("E29").Value = -60
$ws.Range("F29").Value = 22
$ws.Range("G29").Value = 31
$ws.Range("H29").Value = -29.032258064516
$ws.Range("I29").Value = 108
$ws.Range("J29").Value = 139
$ws.Range("K29").Value = -22.302158273381
$ws.Range("L29").Value = -13.6
$ws.Range("M29").Value = -30.769230769230
$ws.Range("N29").Value = -76.774193548387

# ---------------------------------------------------------------------
# Row 30 - Shooting Inc.  (D30 unchanged = 7)
# ---------------------------------------------------------------------
$ws.Range("C30").Value = 4
$ws.Range("E30").Value = -42.857142857142
$ws.Range("F30").Value = 18
$ws.Range("G30").Value = 26
$ws.Range("H30").Value = -30.769230769230
$ws.Range("I30").Value = 93
$ws.Range("J30").Value = 112
$ws.Range("K30").Value = -16.964285714285
$ws.Range("L30").Value = -8.823529411764
$ws.Range("M30").Value = -29.007633587786
$ws.Range("N30").Value = -77.644230769230

# ---------------------------------------------------------------------
# Row 31 - Hate Crimes  (F31 unchanged = 2, J31 unchanged = 11)
# C31 flips from a number (1) to the "***.*" placeholder text.
# ---------------------------------------------------------------------
$ws.Range("C31").Value = "***.*"
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 8
$ws.Range("K31").Value = -27.272727272727
$ws.Range("L31").Value = -20

# ---------------------------------------------------------------------
# Row 33 - Traffic Fatalities  (F33 unchanged = 2, J33 unchanged = 15)
# C33 flips from the "***.*" placeholder text to a number (1).
# ---------------------------------------------------------------------
$ws.Range("C33").Value = 1
$ws.Range("I33").Value = 11
$ws.Range("K33").Value = -26.666666666666
$ws.Range("L33").Value = -26.666666666666
